# feat: add matchsvr gamesvr
#
# Prunes the "全局枚举表" (global enum) sheet down to the core None/1/2/3
# tier rows that matchsvr/gamesvr still need, dropping the higher-tier
# GameType / MatchType / CoinType variants (SNG match type, Short/Aof/Plo*
# game types, Ace/Score/Ticket coin types, etc.).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2: drop the "MatchType SNG" entry (D2); its neighbours (Normal game
# type, Low room type, Coin coin type, Develop env type) stay put. A
# no-op format touch keeps the now-blank cell present (surrounded by
# populated cells) rather than dropped entirely.
$ws.Range("D2").ClearContents()
$ws.Range("D2").Font.Bold = $false

# Row 3: drop the "GameType Short" (A3) and "CoinType Ace" (C3) entries;
# keep "RoomType Middle" (B3) and "EnvType Release" (E3). A3 is the first
# (leading) cell in the row so it drops out entirely; C3 sits between
# populated cells so it stays as a blank placeholder.
$ws.Range("A3").ClearContents()
$ws.Range("C3").ClearContents()
$ws.Range("C3").Font.Bold = $false

# Row 4: drop the "GameType Aof" (A4) and "CoinType Score" (C4) entries;
# keep "RoomType High" (B4).
$ws.Range("A4").ClearContents()
$ws.Range("C4").ClearContents()

# Rows 5-8 only held extra GameType (Plo/Plo5/Plo6/PointRummy) and the
# CoinType "Ticket" entries - remove the rows outright.
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(5).Delete()

# Restore the saved cursor/selection position on the sheet.
[void]$ws.Range("B12").Select()
